# daily auto push: 2026-02-09 14:19 UTC
# A new measurement for 2026/02/09 (月) was logged, inserted as a new row
# right after the existing 2026/02/09 06:00 entry (row 802), pushing all
# subsequent rows down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 802, shifting rows 802:843 down to 803:844
$ws.Rows("802:802").Insert()

# Populate the newly inserted row with the new reading.
# Use a leading apostrophe so Excel keeps the date-looking text as a literal
# string (matching the rest of column A) instead of auto-converting it to a
# real date serial value.
$ws.Cells.Item(802, 1).Value = "'2026/02/09"
$ws.Cells.Item(802, 2).Value = "月"
$ws.Cells.Item(802, 3).Value = 19
$ws.Cells.Item(802, 4).Value = 201

# Reset the style of the new date cell so it doesn't keep an explicit
# "Text" number format applied by the auto-text coercion above; the other
# data rows in this column carry no explicit cell style either.
$ws.Cells.Item(802, 1).Style = "Normal"
